$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update data rows (row 3 is now "1.5 hrs - Core java - Active - From Beginning") ---
$ws.Range("B3").Value = "4/3/2021"
$ws.Range("C3").Value = 1.5
$ws.Range("D3").Value = "Core java"
$ws.Range("E3").Value = "Active"
$ws.Range("F3").Value = " From Beginning"
$ws.Range("D3").Style = "Check Cell"

$ws.Range("B4").Value = "4/3/2021"
$ws.Range("C4").Value = 2.5
$ws.Range("D4").Value = "Program Practice"
$ws.Range("E4").Value = "……."
$ws.Range("F4").ClearContents()

$ws.Range("B5").Value = "4/2/2021"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "Doubt discuss "
$ws.Range("E5").Value = "Cleared"
$ws.Range("F5").Value = "with Asif"

$ws.Range("B6").Value = "4/3/2021"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "Collection"
$ws.Range("E6").Value = "….."
$ws.Range("F6").ClearContents()

$ws.Range("B7").Value = "4/3/2021"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "Collection"
$ws.Range("E7").Value = "Active"
$ws.Range("F7").Value = "not completed"
$ws.Range("H7").Value = " "

# --- Column widths ---
$ws.Range("D1").ColumnWidth = 24.77734375
$ws.Range("E1").ColumnWidth = 21.44140625
$ws.Range("F1").ColumnWidth = 50.77734375

# --- Add Sheet2 (after Sheet1) ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Sheet2"

# --- Re-activate Sheet1 and set selection ---
$ws.Activate()
$ws.Range("F3").Select()
